$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.601.95"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.160.39"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "226.99"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Value = "62.85"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "15.87"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "2.480.79"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "21.74"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").Value = "0.805"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "2.164.00"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "39.577.45"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "0.0₃0856"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "227.83"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  -6.03%  "
$ws.Range("D26").Value = "170.31"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "9.44"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").Value = "1.42"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").Value = "4.57"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "'4.70"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").Value = "0.0618"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +8.26%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "5.08"
$ws.Range("E39").Value = "  +22.62%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "102.54"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "0.0228"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").Value = "1.514.16"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "7.87"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  +27.50%  "
$ws.Range("E51").Value = "  +0.58%  "
